$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E3").Value = 16.577
$ws.Range("E4").Value = 16.431
$ws.Range("D7").Value = -7.193
$ws.Range("C8").Value = -12.765
$ws.Range("C10").Value = -12.658
$ws.Range("E11").Value = 17.074
$ws.Range("C12").Value = -11.207
$ws.Range("D14").Value = -7.878
$ws.Range("E14").Value = 17.329
$ws.Range("D15").Value = -8.317
$ws.Range("C18").Value = -13.806
$ws.Range("D18").Value = -8.529999999999998
$ws.Range("E18").Value = 16.285
$ws.Range("E19").Value = 16.546
$ws.Range("D20").Value = -7.19
$ws.Range("E21").Value = 16.52
$ws.Range("C25").Value = -11.9
$ws.Range("E27").Value = 16.512
$ws.Range("D29").Value = -7.292
$ws.Range("D30").Value = -7.331
$ws.Range("D31").Value = -7.697
$ws.Range("E31").Value = 16.962
$ws.Range("D35").Value = -7.737
$ws.Range("C37").Value = -13.395
$ws.Range("E38").Value = 16.591
$ws.Range("D40").Value = -7.411
$ws.Range("E42").Value = 16.425
$ws.Range("D44").Value = -7.181
$ws.Range("E44").Value = 16.996
$ws.Range("E47").Value = 16.466
$ws.Range("D50").Value = -8.105
$ws.Range("D54").Value = -8.089
$ws.Range("C55").Value = -14.163
$ws.Range("E56").Value = 16.602
$ws.Range("E58").Value = 16.652
$ws.Range("E65").Value = 17.119
$ws.Range("C68").Value = -11.509
$ws.Range("D68").Value = -7.257
$ws.Range("E73").Value = 16.77
$ws.Range("D76").Value = -7.161
$ws.Range("C77").Value = -13.531
$ws.Range("C78").Value = -13.375
$ws.Range("C79").Value = -13.099
$ws.Range("C80").Value = -13.268
$ws.Range("C81").Value = -13.234
$ws.Range("C82").Value = -11.933
$ws.Range("C84").Value = -12.557
$ws.Range("D87").Value = -8.341000000000001
$ws.Range("D88").Value = -8.119
$ws.Range("E90").Value = 16.492
$ws.Range("D92").Value = -7.517
$ws.Range("E92").Value = 16.442
$ws.Range("E94").Value = 17.741
$ws.Range("E95").Value = 17.121
$ws.Range("D96").Value = -7.267
$ws.Range("D98").Value = -8.197999999999999
$ws.Range("C101").Value = -12.573
$ws.Range("D101").Value = -7.678999999999999
$ws.Range("E101").Value = 16.655
$ws.Range("C102").Value = -13.836
$ws.Range("D102").Value = -7.81